# "Añadir configuracion de colegios electorales"
# Replace the "political parties" sample sheet with a "polling place /
# constituency" configuration table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Wipe the old sample data (rows 3:9) and the old header row contents so we
# can lay down the new table from scratch.
$ws.Range("A1:C9").ClearContents()

# New header row: Region | Constitutency | Polling Place
$ws.Range("A1").Value = "Region"
$ws.Range("B1").Value = "Constitutency"
$ws.Range("C1").Value = "Polling Place"

# Row 2 is intentionally left blank (separator row).

# New data rows.
$ws.Range("A3").Value = "Asturias"
$ws.Range("B3").Value = "A"
$ws.Range("C3").Value = 1

$ws.Range("A4").Value = "Madrid"
$ws.Range("B4").Value = "M"
$ws.Range("C4").Value = 2

$ws.Range("A5").Value = "Extremadura"
$ws.Range("B5").Value = "E"
$ws.Range("C5").Value = 3

# Columns A and C go back to the workbook default width; column B widens to
# fit the new "Constitutency" header.
$ws.Columns.Item(1).ColumnWidth = $ws.StandardWidth
$ws.Columns.Item(3).ColumnWidth = $ws.StandardWidth
$ws.Columns.Item(2).ColumnWidth = 14.8

# Selection ends on the new entry table (C3:C5), active cell on the last one.
$ws.Range("C3:C5").Select()
